# MirrorMe workbook edit
# - Split the single "rebuttal" row (row 4) into two rows: the existing
#   rebuttal stays in H4, and a new row 5 is inserted to hold the
#   second rebuttal (previously in G4) in column H.
# - This pushes the "[Scope]" / "[Ttext]" tables below down by one row.
# - Update the active selection to A10 to match the new layout.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remember the text that currently lives in G4 before we start moving rows.
$secondRebuttalText = $ws.Range("G4").Value2

# Insert a new row above row 5. Everything at row 5 and below (including the
# "[Scope]" and "[Ttext]" tables) shifts down by one row.
$ws.Rows.Item(5).Insert()

# The insert copies row 4's formatting into every column of the new row that
# had a value in row 4 (A, D, E, F, G). None of those should hold data in the
# new row, so clear them back out, keeping only H5.
$ws.Range("A5:G5").Clear()

# Clear the old G4 value (it has moved to the new row) and place the second
# rebuttal text into H5 instead.
$ws.Range("G4").Clear()
$ws.Range("H5").Value = $secondRebuttalText

# Match the formatting of the row above (row 4): same row height, wrapped /
# top-aligned cell style.
$ws.Rows.Item(5).RowHeight = $ws.Rows.Item(4).RowHeight
$ws.Range("H5").WrapText = $true
$ws.Range("H5").VerticalAlignment = -4160

# Update dimension / selection to reflect the new layout.
$ws.Range("A10").Select()
